$p = $ppt.ActivePresentation

$oldPathFragment = "integral-private\vignettes\ 2022-09-20"
$newPathFragment = "integral-private\vignettes\C:\Users\kheal579\AppData\Local\Temp\RtmpEzTRKa\callr-scr-94485bb46ca9 2022-08-24"

for ($si = 1; $si -le $p.Slides.Count; $si++) {
    $s = $p.Slides.Item($si)

    for ($i = 1; $i -le $s.Shapes.Count; $i++) {
        $sh = $s.Shapes.Item($i)

        if (-not $sh.HasTextFrame) {
            continue
        }

        $tr = $sh.TextFrame.TextRange

        if ($sh.Name -eq "Title 16") {
            # Paragraph runs today:
            #   1) "Figure 1.\r"                       (bold)
            #   2) "An Empty/A Second Plot in a Template"   (italic)
            #   3) ""                                  (plain)
            #   4) ""                                  (plain, duplicate of 3)
            #
            # Target adds a trailing line break to run 2 (kept italic) and
            # turns the blank "line 3" into its own plain-formatted run
            # holding that extra line break.
            $titleRun = $tr.Runs(2, 1)
            $titleRun.Text = $titleRun.Text + "`r"

            $titleRun = $tr.Runs(2, 1)
            [void]$titleRun.InsertAfter("X")

            $blankRun = $tr.Characters($tr.Length, 1)
            $blankRun.Text = "`r"
            $blankRun.Font.Italic = $false
            $blankRun.Font.Bold = $false
        }
        elseif ($sh.Name -eq "Path Placeholder") {
            if ($tr.Text.Contains($oldPathFragment)) {
                $tr.Text = $tr.Text.Replace($oldPathFragment, $newPathFragment)
            }
        }
    }
}
